$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look numeric must be forced to Text format
# first, otherwise Excel auto-converts them to floating point numbers and
# loses the exact literal representation (e.g. "41.60" -> 41.6).
$textCells = @("D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D17", "D18", "D20", "D21", "D22", "D25", "D26", "D27", "D29", "D30", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values scraped from coinranking.com
$ws.Range("D2").Value = '28.427.02'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '1.825.39'
$ws.Range("E3").Value = '  -0.64%  '
$ws.Range("E4").Value = '  +0.40%  '
$ws.Range("D5").Value = '315.13'
$ws.Range("E5").Value = '  -1.23%  '
$ws.Range("E6").Value = '  +0.31%  '
$ws.Range("D7").Value = '0.5114'
$ws.Range("E7").Value = '  -3.88%  '
$ws.Range("D8").Value = '0.3919'
$ws.Range("E8").Value = '  -1.67%  '
$ws.Range("D9").Value = '0.07648'
$ws.Range("E9").Value = '  +0.87%  '
$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").Value = '41.60'
$ws.Range("E10").Value = '  -0.50%  '
$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D11").Value = '1.107'
$ws.Range("E11").Value = '  +0.03%  '
$ws.Range("D12").Value = '20.98'
$ws.Range("E12").Value = '  +0.93%  '
$ws.Range("D13").Value = '6.271'
$ws.Range("E13").Value = '  -0.76%  '
$ws.Range("D14").Value = '1.003'
$ws.Range("E14").Value = '  +0.44%  '
$ws.Range("D15").Value = '7.511'
$ws.Range("E15").Value = '  -1.44%  '
$ws.Range("D16").Value = '1.821.98'
$ws.Range("E16").Value = '  -0.27%  '
$ws.Range("D17").Value = '93.41'
$ws.Range("E17").Value = '  +3.89%  '
$ws.Range("D18").Value = '0.00001097'
$ws.Range("E18").Value = '  +2.22%  '
$ws.Range("E19").Value = '  +1.14%  '
$ws.Range("D20").Value = '17.63'
$ws.Range("E20").Value = '  -0.14%  '
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.26%  '
$ws.Range("D22").Value = '6.143'
$ws.Range("E22").Value = '  +1.28%  '
$ws.Range("D23").Value = '28.464.01'
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("E24").Value = '  -0.47%  '
$ws.Range("D25").Value = '2.265'
$ws.Range("E25").Value = '  +7.75%  '
$ws.Range("D26").Value = '20.67'
$ws.Range("E26").Value = '  +0.33%  '
$ws.Range("D27").Value = '156.78'
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("D28").Value = '2.035.04'
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("D29").Value = '2.379'
$ws.Range("E29").Value = '  -2.42%  '
$ws.Range("D30").Value = '124.22'
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("E31").Value = '  -0.72%  '
$ws.Range("E32").Value = '  -1.29%  '
$ws.Range("D33").Value = '5.632'
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("D34").Value = '3.656'
$ws.Range("E34").Value = '  -0.95%  '
$ws.Range("D35").Value = '0.07024'
$ws.Range("E35").Value = '  -4.25%  '
$ws.Range("D36").Value = '0.2204'
$ws.Range("E36").Value = '  -2.08%  '
$ws.Range("D37").Value = '8.881'
$ws.Range("E37").Value = '  +0.17%  '
$ws.Range("D38").Value = '0.02315'
$ws.Range("E38").Value = '  -0.91%  '
$ws.Range("E39").Value = '  -1.91%  '
$ws.Range("D40").Value = '0.6238'
$ws.Range("E40").Value = '  -0.70%  '
$ws.Range("D41").Value = '11.20'
$ws.Range("E41").Value = '  -1.47%  '
$ws.Range("D42").Value = '1.172'
$ws.Range("E42").Value = '  -2.09%  '
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.21%  '
$ws.Range("D44").Value = '1.390'
$ws.Range("E44").Value = '  -1.84%  '
$ws.Range("D45").Value = '13.40'
$ws.Range("E45").Value = '  -0.53%  '
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").Value = '3.711'
$ws.Range("E46").Value = '  +0.15%  '
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '0.5881'
$ws.Range("E47").Value = '  +0.94%  '
$ws.Range("D48").Value = '125.19'
$ws.Range("E48").Value = '  -0.51%  '
$ws.Range("D49").Value = '1.971'
$ws.Range("E49").Value = '  -0.19%  '
$ws.Range("D50").Value = '1.194'
$ws.Range("E50").Value = '  +0.14%  '
$ws.Range("D51").Value = '0.06924'
$ws.Range("E51").Value = '  +0.26%  '
